$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.173.86'
$ws.Range("E2").Value = '  -0.72%  '
$ws.Range("D3").Value = '1.906.76'
$ws.Range("E3").Value = '  -1.63%  '
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7299'
$ws.Range("E5").Value = '  -5.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.63'
$ws.Range("E6").Value = '  -2.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3121'
$ws.Range("E8").Value = '  -2.98%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '26.42'
$ws.Range("E9").Value = '  -5.93%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06887'
$ws.Range("E10").Value = '  -3.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7762'
$ws.Range("E11").Value = '  -1.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07944'
$ws.Range("E12").Value = '  -1.09%  '
$ws.Range("D13").Value = '1.886.16'
$ws.Range("E13").Value = '  -2.77%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.255'
$ws.Range("E14").Value = '  -2.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.25'
$ws.Range("E15").Value = '  -4.07%  '
$ws.Range("D16").Value = '30.118.73'
$ws.Range("E16").Value = '  -0.91%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.13'
$ws.Range("E17").Value = '  -3.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.826'
$ws.Range("E18").Value = '  -0.28%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '239.42'
$ws.Range("E19").Value = '  -6.56%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007758'
$ws.Range("E20").Value = '  -3.40%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  +0.26%  '
$ws.Range("D22").Value = '2.140.38'
$ws.Range("E22").Value = '  -2.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.002'
$ws.Range("E23").Value = '  +0.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.951'
$ws.Range("E24").Value = '  +2.69%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.359'
$ws.Range("E25").Value = '  -2.84%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.58'
$ws.Range("E26").Value = '  +0.29%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.08'
$ws.Range("E27").Value = '  -0.43%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1272'
$ws.Range("E28").Value = '  -5.51%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.048'
$ws.Range("E29").Value = '  -11.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.357'
$ws.Range("E30").Value = '  -0.61%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.540'
$ws.Range("E31").Value = '  +0.91%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.292'
$ws.Range("E32").Value = '  -3.37%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.079'
$ws.Range("E33").Value = '  -1.82%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05126'
$ws.Range("E34").Value = '  -1.50%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.287'
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7378'
$ws.Range("E36").Value = '  -2.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.751'
$ws.Range("E37").Value = '  -0.79%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01929'
$ws.Range("E38").Value = '  -2.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.783'
$ws.Range("E39").Value = '  -1.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.349'
$ws.Range("E40").Value = '  -2.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.66'
$ws.Range("E41").Value = '  -5.64%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4428'
$ws.Range("E42").Value = '  -2.31%  '
$ws.Range("E43").Value = '  -2.78%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8352'
$ws.Range("E45").Value = '  -0.17%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.04'
$ws.Range("E46").Value = '  -0.46%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.559'
$ws.Range("E47").Value = '  +0.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.690'
$ws.Range("E48").Value = '  -1.14%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '37.53'
$ws.Range("E49").Value = '  +0.07%  '
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '939.85'
$ws.Range("E50").Value = '  -4.61%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1183'
$ws.Range("E51").Value = '  -0.42%  '
